# Generate Report for Handoff
# b.md has been handed off again (new handback xliff generated for zh-cn and
# de-de, but the source has drifted so the status reverts to "Ready for
# handoff" and an error detail is recorded).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9c6df2c8901979d691951445d8a500835290073/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a6d55190db70616c41220e64bec2a64996852a6/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: refresh the b.md row's status + generate-date columns
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-03 10:41:04"

# ---------------------------------------------------------------------
# zh-cn sheet: new handoff xliff for b.md, status reset, error populated
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# "False" is stored as plain text (like the rest of the "Content
# Duplicate" column), not a native boolean, so copy it from a sibling
# cell that already holds that text instead of typing it (typing a bare
# True/False literal gets auto-coerced to a real boolean by Excel).
$wsZhCn.Range("O2").Copy()
$wsZhCn.Range("F3").PasteSpecial(-4163)
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-03 10:40:58"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666667

# ---------------------------------------------------------------------
# de-de sheet: new handoff xliff for b.md, status reset, error populated
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("O2").Copy()
$wsDeDe.Range("F3").PasteSpecial(-4163)
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-03 10:41:04"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666667

$excel.CutCopyMode = $false
